# Add a new "Emp" worksheet after the existing "VerifyLoginLogout" sheet,
# populate it with employee first/last name sample data, and make it the
# active sheet with B4 selected (matching the authored workbook state).

$wb = $excel.ActiveWorkbook

$firstSheet = $wb.ActiveSheet

# Insert the new sheet right after the current (first/only) sheet so it
# lands at the end of the tab strip, like Excel's own "Insert Sheet" button.
$empSheet = $wb.Worksheets.Add($null, $firstSheet)
$empSheet.Name = "Emp"

# Header row
$empSheet.Range("A1").Value = "FN"
$empSheet.Range("B1").Value = "LN"

# Sample rows
$empSheet.Range("A2").Value = "a"
$empSheet.Range("B2").Value = "a"
$empSheet.Range("A3").Value = "b"
$empSheet.Range("B3").Value = "b"
$empSheet.Range("A4").Value = "c"
$empSheet.Range("B4").Value = "c"

# Leave the new sheet active/selected, with B4 as the last selected cell.
$empSheet.Activate() | Out-Null
$empSheet.Range("B4").Select() | Out-Null
